$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains text formatting so numeric-looking
# strings (e.g. "1.00", "0.615") are not silently converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "38.752.70"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "2.087.79"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "228.44"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "2.398.08"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E13").Value = "  +3.96%  "
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").Value = "0.796"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("D16").Value = "5.48"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "2.090.21"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "38.682.41"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").Value = "71.50"
$ws.Range("E19").Value = "  +2.87%  "
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "227.20"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "170.88"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "9.52"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").Value = "0.139"
$ws.Range("E28").Value = "  +9.16%  "
$ws.Range("D29").Value = "1.46"
$ws.Range("E29").Value = "  +12.95%  "
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "2.39"
$ws.Range("E32").Value = "  +6.05%  "
$ws.Range("D33").Value = "4.49"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("D35").Value = "0.0608"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").Value = "6.47"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "18.03"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").Value = "1.542.01"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").Value = "100.65"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "0.0921"
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("D46").Value = "7.69"
$ws.Range("E46").Value = "  +8.29%  "
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "4.13"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "2.286.77"
$ws.Range("E51").Value = "  +2.13%  "
